$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 4103.1377
$ws.Range("I62").Value = 1738.0952
$ws.Range("J62").Value = 10311.375
$ws.Range("K62").Value = 1738.0952
$ws.Range("L62").Value = 10311.375
$ws.Range("M62").Value = -1114.0952
$ws.Range("N62").Value = -11559.375
$ws.Range("H65").Value = 4103.1377
$ws.Range("I65").Value = 1738.0952
$ws.Range("J65").Value = 10311.375
$ws.Range("K65").Value = 8690.476000000001
$ws.Range("L65").Value = 51556.875
$ws.Range("M65").Value = -5570.476000000001
$ws.Range("N65").Value = -57796.875
$ws.Range("H74").Value = 3876.4119
$ws.Range("I74").Value = 3766.5557
$ws.Range("J74").Value = 4000
$ws.Range("K74").Value = 3766.5557
$ws.Range("L74").Value = 4000
$ws.Range("M74").Value = -2830.5557
$ws.Range("N74").Value = -5872
$ws.Range("H77").Value = 3876.4119
$ws.Range("I77").Value = 3766.5557
$ws.Range("J77").Value = 4000
$ws.Range("K77").Value = 18832.7785
$ws.Range("L77").Value = 20000
$ws.Range("M77").Value = -14152.7785
$ws.Range("N77").Value = -29360
$ws.Range("H111").Value = 1642.7142
$ws.Range("J111").Value = 1479.8
$ws.Range("L111").Value = 4439.4
$ws.Range("N111").Value = -10573.4
$ws.Range("H137").Value = 142859790
$ws.Range("I137").Value = 166668100
$ws.Range("J137").Value = 10003
$ws.Range("K137").Value = 500004300
$ws.Range("L137").Value = 30009
$ws.Range("M137").Value = -500001750
$ws.Range("N137").Value = -35109
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 787.2
$ws.Range("I45").Value = 735.25
$ws.Range("K45").Value = 735.25
$ws.Range("M45").Value = -358.25
$ws.Range("H61").Value = 2792.875
$ws.Range("I61").Value = 1652.4
$ws.Range("J61").Value = 4693.6665
$ws.Range("K61").Value = 1652.4
$ws.Range("L61").Value = 4693.6665
$ws.Range("M61").Value = -1440.4
$ws.Range("N61").Value = -5117.6665
$ws.Range("H110").Value = 525.1739
$ws.Range("I110").Value = 498.95
$ws.Range("K110").Value = 498.95
$ws.Range("M110").Value = 1546.05
$ws.Range("H122").Value = 1673.4722
$ws.Range("I122").Value = 1453.9524
$ws.Range("J122").Value = 1980.8
$ws.Range("K122").Value = 4361.857199999999
$ws.Range("L122").Value = 5942.4
$ws.Range("M122").Value = -1911.857199999999
$ws.Range("N122").Value = -10842.4
$ws.Range("H136").Value = 2792.875
$ws.Range("I136").Value = 1652.4
$ws.Range("J136").Value = 4693.6665
$ws.Range("K136").Value = 4957.200000000001
$ws.Range("L136").Value = 14080.9995
$ws.Range("M136").Value = -2407.200000000001
$ws.Range("N136").Value = -19180.9995
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1463.6
$ws.Range("I20").Value = 1514.3334
$ws.Range("J20").Value = 1387.5
$ws.Range("K20").Value = 1514.3334
$ws.Range("L20").Value = 1387.5
$ws.Range("M20").Value = -1267.3334
$ws.Range("N20").Value = -1881.5
$ws.Range("H94").Value = 1188.2858
$ws.Range("I94").Value = 1421.3334
$ws.Range("J94").Value = 768.8
$ws.Range("K94").Value = 1421.3334
$ws.Range("L94").Value = 768.8
$ws.Range("M94").Value = -970.3334
$ws.Range("N94").Value = -1670.8
$ws.Range("H105").Value = 3710.8696
$ws.Range("I105").Value = 3566.6667
$ws.Range("J105").Value = 3981.25
$ws.Range("K105").Value = 3566.6667
$ws.Range("L105").Value = 3981.25
$ws.Range("M105").Value = -1819.6667
$ws.Range("N105").Value = -7475.25
$ws.Range("H107").Value = 869.44446
$ws.Range("I107").Value = 691.4286
$ws.Range("J107").Value = 1492.5
$ws.Range("K107").Value = 691.4286
$ws.Range("L107").Value = 1492.5
$ws.Range("M107").Value = 1228.5714
$ws.Range("N107").Value = -5332.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 663.3
$ws.Range("I22").Value = 403.66666
$ws.Range("J22").Value = 3000
$ws.Range("K22").Value = 403.66666
$ws.Range("L22").Value = 3000
$ws.Range("M22").Value = -53.66665999999998
$ws.Range("N22").Value = -3700
$ws.Range("H94").Value = 1555.421
$ws.Range("I94").Value = 977.875
$ws.Range("K94").Value = 977.875
$ws.Range("M94").Value = -526.875
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H82").Value = 2868
$ws.Range("I82").Value = 406.5
$ws.Range("J82").Value = 3571.2856
$ws.Range("K82").Value = 1219.5
$ws.Range("L82").Value = 10713.8568
$ws.Range("M82").Value = -813.5
$ws.Range("N82").Value = -11525.8568
$ws.Range("H85").Value = 2868
$ws.Range("I85").Value = 406.5
$ws.Range("J85").Value = 3571.2856
$ws.Range("K85").Value = 1219.5
$ws.Range("L85").Value = 10713.8568
$ws.Range("M85").Value = 184.5
$ws.Range("N85").Value = -13521.8568
$ws.Range("H127").Value = 2081.4443
$ws.Range("J127").Value = 2081.4443
$ws.Range("L127").Value = 6244.3329
$ws.Range("N127").Value = -16164.3329
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1724.75
$ws.Range("I97").Value = 1256.8572
$ws.Range("J97").Value = 5000
$ws.Range("K97").Value = 1256.8572
$ws.Range("L97").Value = 5000
$ws.Range("M97").Value = -760.8571999999999
$ws.Range("N97").Value = -5992
$ws.Range("H102").Value = 2029
$ws.Range("I102").Value = 1742.1666
$ws.Range("K102").Value = 1742.1666
$ws.Range("M102").Value = -120.1666
$ws.Range("H122").Value = 1853336.4
$ws.Range("I122").Value = 2778879.5
$ws.Range("K122").Value = 8336638.5
$ws.Range("M122").Value = -8334188.5
$ws.Range("H132").Value = 3717
$ws.Range("I132").Value = 3613
$ws.Range("K132").Value = 10839
$ws.Range("M132").Value = -8309
$ws.Range("H139").Value = 39550.332
$ws.Range("J139").Value = 39550.332
$ws.Range("L139").Value = 39550.332
$ws.Range("N139").Value = -49830.332
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2798.3044
$ws.Range("I7").Value = 1407.625
$ws.Range("K7").Value = 1407.625
$ws.Range("M7").Value = -1295.625
$ws.Range("H40").Value = 2779.4075
$ws.Range("I40").Value = 1683.4286
$ws.Range("J40").Value = 3163
$ws.Range("K40").Value = 1683.4286
$ws.Range("L40").Value = 3163
$ws.Range("M40").Value = -1547.4286
$ws.Range("N40").Value = -3435
$ws.Range("H61").Value = 3400.4
$ws.Range("I61").Value = 3386
$ws.Range("J61").Value = 3447.7144
$ws.Range("K61").Value = 3386
$ws.Range("L61").Value = 3447.7144
$ws.Range("M61").Value = -3184
$ws.Range("N61").Value = -3851.7144
$ws.Range("H93").Value = 927.6842
$ws.Range("J93").Value = 1330.4
$ws.Range("L93").Value = 1330.4
$ws.Range("N93").Value = -3826.4
$ws.Range("H113").Value = 3400.4
$ws.Range("I113").Value = 3386
$ws.Range("J113").Value = 3447.7144
$ws.Range("K113").Value = 3386
$ws.Range("L113").Value = 3447.7144
$ws.Range("M113").Value = -1216
$ws.Range("N113").Value = -7787.7144
$ws.Range("H122").Value = 2888.8462
$ws.Range("I122").Value = 1807.0714
$ws.Range("J122").Value = 3494.64
$ws.Range("K122").Value = 5421.2142
$ws.Range("L122").Value = 10483.92
$ws.Range("M122").Value = -2971.2142
$ws.Range("N122").Value = -15383.92
$ws.Range("H126").Value = 2798.3044
$ws.Range("I126").Value = 1407.625
$ws.Range("K126").Value = 4222.875
$ws.Range("M126").Value = -1752.875
$ws.Range("H132").Value = 3990.383
$ws.Range("I132").Value = 3366.6365
$ws.Range("J132").Value = 4539.28
$ws.Range("K132").Value = 10099.9095
$ws.Range("L132").Value = 13617.84
$ws.Range("M132").Value = -7569.9095
$ws.Range("N132").Value = -18677.84
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 11127819
$ws.Range("I62").Value = 14304208
$ws.Range("J62").Value = 10455.5
$ws.Range("K62").Value = 14304208
$ws.Range("L62").Value = 10455.5
$ws.Range("M62").Value = -14303584
$ws.Range("N62").Value = -11703.5
$ws.Range("H65").Value = 11127819
$ws.Range("I65").Value = 14304208
$ws.Range("J65").Value = 10455.5
$ws.Range("K65").Value = 71521040
$ws.Range("L65").Value = 52277.5
$ws.Range("M65").Value = -71517920
$ws.Range("N65").Value = -58517.5
$ws.Range("H92").Value = 26758.166
$ws.Range("J92").Value = 26758.166
$ws.Range("L92").Value = 26758.166
$ws.Range("N92").Value = -31750.166
$ws.Range("H122").Value = 24031.268
$ws.Range("I122").Value = 38503
$ws.Range("J122").Value = 2323.6667
$ws.Range("K122").Value = 115509
$ws.Range("L122").Value = 6971.000100000001
$ws.Range("M122").Value = -113059
$ws.Range("N122").Value = -11871.0001
$ws.Range("H126").Value = 45210.566
$ws.Range("I126").Value = 49368.715
$ws.Range("J126").Value = 1550
$ws.Range("K126").Value = 148106.145
$ws.Range("L126").Value = 4650
$ws.Range("M126").Value = -145636.145
$ws.Range("N126").Value = -9590
$ws.Range("H132").Value = 20005038
$ws.Range("I132").Value = 45462656
$ws.Range("J132").Value = 2625.7144
$ws.Range("K132").Value = 136387968
$ws.Range("L132").Value = 7877.1432
$ws.Range("M132").Value = -136385438
$ws.Range("N132").Value = -12937.1432
$ws.Range("H136").Value = 10786891
$ws.Range("J136").Value = 2713.6155
$ws.Range("L136").Value = 8140.8465
$ws.Range("N136").Value = -13240.8465
